$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D keeps its existing text formatting; only cells whose new
# value would otherwise be auto-converted to a number get an explicit
# text number-format applied first so the stored value stays a string.

$ws.Range("D2").Value = "43.988.96"
$ws.Range("E2").Value = "  +2.52%  "

$ws.Range("D3").Value = "2.255.64"
$ws.Range("E3").Value = "  +1.79%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.58"
$ws.Range("E5").Value = "  +0.14%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.93"
$ws.Range("E6").Value = "  +3.82%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.577"
$ws.Range("E7").Value = "  -0.14%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.554"
$ws.Range("E9").Value = "  -0.06%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.46"
$ws.Range("E10").Value = "  +2.50%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0832"
$ws.Range("E11").Value = "  +1.08%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.62"
$ws.Range("E12").Value = "  +0.51%  "

$ws.Range("E13").Value = "  -1.14%  "

$ws.Range("D14").Value = "2.599.91"
$ws.Range("E14").Value = "  +1.67%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.858"
$ws.Range("E15").Value = "  -0.02%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.46"
$ws.Range("E16").Value = "  +1.14%  "

$ws.Range("D17").Value = "2.259.99"
$ws.Range("E17").Value = "  +1.54%  "

$ws.Range("D18").Value = "43.889.55"
$ws.Range("E18").Value = "  +2.33%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.43"
$ws.Range("E19").Value = "  -4.28%  "

$ws.Range("D20").Value = "0.0₃0986"
$ws.Range("E20").Value = "  +2.81%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.53"
$ws.Range("E21").Value = "  +0.35%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.84"
$ws.Range("E22").Value = "  +1.74%  "

$ws.Range("E23").Value = "  -0.63%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "235.33"
$ws.Range("E24").Value = "  +0.17%  "

$ws.Range("E25").Value = "  -1.62%  "

$ws.Range("E26").Value = "  +0.22%  "

$ws.Range("E27").Value = "  +3.39%  "

$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.20"
$ws.Range("E28").Value = "  -1.81%  "

$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "37.92"
$ws.Range("E29").Value = "  +7.04%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.21"
$ws.Range("E30").Value = "  -1.11%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "160.88"
$ws.Range("E31").Value = "  +6.24%  "

$ws.Range("E32").Value = "  +0.39%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0852"
$ws.Range("E33").Value = "  -0.97%  "

$ws.Range("E34").Value = "  +1.29%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.116"
$ws.Range("E35").Value = "  +12.06%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.94"
$ws.Range("E36").Value = "  +2.46%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.03"
$ws.Range("E37").Value = "  -3.47%  "

$ws.Range("E38").Value = "  -1.37%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.68"
$ws.Range("E39").Value = "  +23.41%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.73"
$ws.Range("E40").Value = "  +2.32%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.22"
$ws.Range("E41").Value = "  -3.84%  "

$ws.Range("E42").Value = "  -0.93%  "

$ws.Range("D44").Value = "1.797.01"
$ws.Range("E44").Value = "  +4.02%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "75.52"
$ws.Range("E45").Value = "  +2.64%  "

$ws.Range("E46").Value = "  -1.59%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "82.81"
$ws.Range("E47").Value = "  -0.81%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.22"
$ws.Range("E48").Value = "  -0.76%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "104.96"
$ws.Range("E49").Value = "  +2.39%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.68"
$ws.Range("E50").Value = "  +8.95%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "58.25"
$ws.Range("E51").Value = "  +1.72%  "
